$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.33483695883526
$ws.Range("B3").Value = 0.5292499577636943
$ws.Range("B4").Value = 16.78528138550232

$ws.Range("A5:B5").Delete()
